$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add Sheet3 at the end of the workbook (after the current last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# Header row (no / name / tahun)
$ws3.Range("A1").Value = "no"
$ws3.Range("B1").Value = "name"
$ws3.Range("C1").Value = "tahun"

# Year row
$ws3.Range("C2").Value = 2020
$ws3.Range("D2").Value = 2021
$ws3.Range("E2").Value = 2022
$ws3.Range("F2").Value = 2023

# Data rows: no, name, year1..year4
$data = @(
    @(1, "John", 23, 24, 25, 26),
    @(2, "Doe", 24, 25, 26, 27),
    @(3, "James", 24, 25, 26, 27),
    @(4, "Mark", 43, 44, 45, 46),
    @(5, "Klare", 54, 55, 56, 57),
    @(6, "Hary", 32, 33, 34, 35),
    @(7, "Ozil", 12, 13, 14, 15)
)

$row = 3
foreach ($r in $data) {
    $ws3.Range("A$row").Value = $r[0]
    $ws3.Range("B$row").Value = $r[1]
    $ws3.Range("C$row").Value = $r[2]
    $ws3.Range("D$row").Value = $r[3]
    $ws3.Range("E$row").Value = $r[4]
    $ws3.Range("F$row").Value = $r[5]
    $row++
}

# Merge the header cells: no/name span two rows, tahun spans the four year columns
$ws3.Range("C1:F1").Merge()
$ws3.Range("A1:A2").Merge()
$ws3.Range("B1:B2").Merge()

# Match existing formatting used across the other sheets (reuse cell style
# rather than minting a new one). Only touch the cells that actually hold
# data, so we don't materialize extra empty-but-styled cells.
$ws1.Range("A1").Copy()
$ws3.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("C2:F2").PasteSpecial(-4122)
$ws3.Range("A3:F9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
